# إضافة حدث جديد في Card21 by admin at 2025-12-08 08:44:54
#
# Card21's service-log table currently ends at row 20 (header in row 1,
# data rows 2..20). This adds a new service event as row 21 and fills in
# the (previously blank) B20:K20 measurement columns with "nan" so the
# sheet matches the rest of the log's placeholder convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card21")

# Row 20 was left with blank placeholder cells in columns B..K; backfill
# them with the "nan" placeholder used throughout the rest of the sheet.
$blankCols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")
foreach ($col in $blankCols) {
    $ws.Range($col + "20").Value = "nan"
}

# New service event -> row 21.
$ws.Range("A21").Value = "21"
$ws.Range("L21").Value = "31\8\2025"
$ws.Range("M21").Value = "796 t"
$ws.Range("N21").Value = "تم تغيير الجرائد الخلفيه (1_5_8)"
$ws.Range("O21").Value = "الخبير"
